$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

# Update the text of C8: the use case's actor action description was reworded.
$ws.Range("C8").Value = "Indica que pretende ver todos os serviços agendados"

# Move the active selection to C9 (as recorded in the saved view state).
$ws.Activate()
$ws.Range("C9").Select()
